$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 57
$ws.Range("I2").Value = 106
$ws.Range("J2").Value = 436
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 121
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 76
$ws.Range("P2").Value = 3
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 58
$ws.Range("T2").Value = 86
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 678
$ws.Range("X2").Value = 708
$ws.Range("Y2").Value = 2
$ws.Range("AA2").Value = 7
